# Generate Report for Handback
# Populates the "Latest Target File" (F) and "Latest Handback File" (G) columns
# for both the zh-cn and de-de sheets, marks both localized files as handed
# back (in sync with en-US), and stamps the handback datetime.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

$statusText = "Handed back: in sync with en-US"

# --- Status column (C on zh-cn/de-de; B & C on Overview) ---
# All of these cells share the same underlying string ("Ready for handoff")
# in the source workbook, so every reference needs to move to the new text
# together for the shared string to collapse back into a single entry.
$wsOverview.Range("B2").Value = $statusText
$wsOverview.Range("C2").Value = $statusText
$wsOverview.Range("B3").Value = $statusText
$wsOverview.Range("C3").Value = $statusText
$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText
$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

# --- Latest Handback DateTime column (H) ---
$wsZh.Range("H2").Value = "2016-03-30 10:08:22"
$wsZh.Range("H3").Value = "2016-03-30 10:08:22"
$wsDe.Range("H2").Value = "2016-03-30 10:08:40"
$wsDe.Range("H3").Value = "2016-03-30 10:08:40"

# --- Latest Target File (F) / Latest Handback File (G) hyperlinks ---

# zh-cn, row 2 (11d03539-...)
$wsZh.Hyperlinks.Add(
    $wsZh.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/b07f9d48f6306d651e348e6b04cb73123700c202/e2e/11d03539-d425-4e32-b99a-31afb4d274be.md",
    [System.Reflection.Missing]::Value,
    [System.Reflection.Missing]::Value,
    "11d03539-d425-4e32-b99a-31afb4d274be.md"
)
$wsZh.Hyperlinks.Add(
    $wsZh.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/494e4fd141e9353437f930c391dc332467da0f07/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/11d03539-d425-4e32-b99a-31afb4d274be.95a75a42193a4a2413bf33e37c089ef2f6232534.zh-cn.xlf",
    [System.Reflection.Missing]::Value,
    [System.Reflection.Missing]::Value,
    "11d03539-d425-4e32-b99a-31afb4d274be.95a75a42193a4a2413bf33e37c089ef2f6232534.zh-cn.xlf"
)

# zh-cn, row 3 (f42e4ba0-...)
$wsZh.Hyperlinks.Add(
    $wsZh.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/b07f9d48f6306d651e348e6b04cb73123700c202/e2e/f42e4ba0-6fd6-4fac-a191-f4cba99ef29c.md",
    [System.Reflection.Missing]::Value,
    [System.Reflection.Missing]::Value,
    "f42e4ba0-6fd6-4fac-a191-f4cba99ef29c.md"
)
$wsZh.Hyperlinks.Add(
    $wsZh.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/494e4fd141e9353437f930c391dc332467da0f07/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/f42e4ba0-6fd6-4fac-a191-f4cba99ef29c.e3cbd3baaa7ed8b866f24df35a3ceef84ff71312.zh-cn.xlf",
    [System.Reflection.Missing]::Value,
    [System.Reflection.Missing]::Value,
    "f42e4ba0-6fd6-4fac-a191-f4cba99ef29c.e3cbd3baaa7ed8b866f24df35a3ceef84ff71312.zh-cn.xlf"
)

# de-de, row 2 (11d03539-...)
$wsDe.Hyperlinks.Add(
    $wsDe.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/b07f9d48f6306d651e348e6b04cb73123700c202/e2e/11d03539-d425-4e32-b99a-31afb4d274be.md",
    [System.Reflection.Missing]::Value,
    [System.Reflection.Missing]::Value,
    "11d03539-d425-4e32-b99a-31afb4d274be.md"
)
$wsDe.Hyperlinks.Add(
    $wsDe.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ac5fc2f4e2e700e5dbb0937fc8280135194291b1/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/11d03539-d425-4e32-b99a-31afb4d274be.95a75a42193a4a2413bf33e37c089ef2f6232534.de-de.xlf",
    [System.Reflection.Missing]::Value,
    [System.Reflection.Missing]::Value,
    "11d03539-d425-4e32-b99a-31afb4d274be.95a75a42193a4a2413bf33e37c089ef2f6232534.de-de.xlf"
)

# de-de, row 3 (f42e4ba0-...)
$wsDe.Hyperlinks.Add(
    $wsDe.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/b07f9d48f6306d651e348e6b04cb73123700c202/e2e/f42e4ba0-6fd6-4fac-a191-f4cba99ef29c.md",
    [System.Reflection.Missing]::Value,
    [System.Reflection.Missing]::Value,
    "f42e4ba0-6fd6-4fac-a191-f4cba99ef29c.md"
)
$wsDe.Hyperlinks.Add(
    $wsDe.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ac5fc2f4e2e700e5dbb0937fc8280135194291b1/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/f42e4ba0-6fd6-4fac-a191-f4cba99ef29c.e3cbd3baaa7ed8b866f24df35a3ceef84ff71312.de-de.xlf",
    [System.Reflection.Missing]::Value,
    [System.Reflection.Missing]::Value,
    "f42e4ba0-6fd6-4fac-a191-f4cba99ef29c.e3cbd3baaa7ed8b866f24df35a3ceef84ff71312.de-de.xlf"
)
